# Commit: "Add files via upload" — 686-huaweicloud-device_info
# The device-info sheet's PORT columns (G = lan/internal port, H = wan/external
# port) are normalized from specific numeric port numbers to the literal text
# "all" for the rows that describe an "open to all ports" configuration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(11, 12, 13, 14, 15, 18, 25, 26, 27)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "all"   # column G
    $ws.Cells.Item($r, 8).Value = "all"   # column H
}

# Restore the view to the top of the sheet and leave the selection on the
# first empty row below the data table (matches the saved workbook state).
[void]$ws.Range("A1").Select()
[void]$ws.Range("I31").Select()
